$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.168.21"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -5.76%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.185.44"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -7.08%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.34%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.77"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.66%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.52"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -7.39%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.586"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -4.75%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.42%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.178.17"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -7.31%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.592"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -6.95%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "51.64"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -12.82%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.129"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -5.54%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000248"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.34%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.78"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -7.07%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.737.19"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -5.24%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.212.89"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -5.71%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.113"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -8.25%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.315.15"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.93%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.98"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.31%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.79"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -5.16%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.946"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.42%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "359.87"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -5.32%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.69"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.26%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.29"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -5.30%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.82"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.65%  "

# Row 26
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.11"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.53%  "

# Row 27
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.82"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.55%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.57"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.62%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.02"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -6.97%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.04"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -6.85%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "642.80"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -6.54%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.87"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -7.47%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.29"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -9.13%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.96"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.37%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.25%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "56.52"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -8.66%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.01%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.92"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.65%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.369"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -5.79%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.36%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0696"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +8.82%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.855.15"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.60%  "

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -7.74%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.46"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.27%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.61"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.12%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0384"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.86%  "

# Row 47
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.75"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.84%  "

# Row 48
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.52"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -10.55%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.92"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.42%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.29"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.92%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.121"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.65%  "
